$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update education table contents -------------------------------------
# The "PhD - Psychology" (Stirling) and "MSc in Evolutionary Psychology"
# (Liverpool) entries swap row order (PhD now listed first, row 2; MSc
# moves to row 5) and both get expanded descriptive text.
#
# Write row 5 (MSc / Liverpool) first, then row 2 (PhD / Stirling), so the
# two brand-new shared-string entries land in the same append order as the
# canonical file (MSc string appended before the PhD string).

# Row 5: MSc in Evolutionary Psychology -> University of Liverpool
$ws.Range("A5").Value = "MSc in Evolutionary Psychology  \textit{(School of Biological Sciences)}"
$ws.Range("C5").Value = "\href{https://www.liverpool.ac.uk/}{University of Liverpool}"
$ws.Range("D5").Value = " Liverpool, Reino Unido"
$ws.Range("E5").Value = "Supervisor: \href{https://www.scraigroberts.com/}{Prof. S. Craig Roberts}"

# Row 2: PhD - Psychology -> University of Stirling
$ws.Range("A2").Value = "PhD - Psychology \textit{(\href{https://www.stir.ac.uk/about/faculties/natural-sciences/our-research/research-groups/behaviour-and-evolution-research-group/}{Behaviour and Evolution Research Group}, Faculty of Natural Sciences)}"
$ws.Range("C2").Value = "\href{https://www.stir.ac.uk/}{University of Stirling}"
$ws.Range("D2").Value = "Stirling, Reino Unido"
$ws.Range("E2").Value = "Tesis: \href{https://dspace.stir.ac.uk/handle/1893/21102}{\textbf{\textit{Contextual musicality: vocal modulation and its perception in human social interaction}}}"

# --- Row 2 grew a lot of wrapped text in column A, so it now needs to be
# taller (60pt -> 120pt). -----------------------------------------------
$ws.Rows.Item(2).RowHeight = 120

# --- Move the active selection to A2 (was E10). ---------------------------
$ws.Range("A2").Select() | Out-Null
